$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Update cited results (citiranost) for 2024 row - z_clanki count updated
$ws.Range("B8").Value = 904

# Update the active cell selection to reflect the latest manual selection
$ws.Range("I17").Select()
